$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 50-51 (changed cells only, per diff)
$ws.Cells.Item(50, 22).Value = 62547500
$ws.Cells.Item(50, 23).Value = 19880100
$ws.Cells.Item(50, 24).Value = 24253600
$ws.Cells.Item(50, 25).Value = 69588200
$ws.Cells.Item(50, 26).Value = 142215700

# Row 51: B-Z changed (A unchanged)
$ws.Cells.Item(51, 2).Value = 209.6799926757812
$ws.Cells.Item(51, 3).Value = 164.7299957275391
$ws.Cells.Item(51, 4).Value = 378.7699890136719
$ws.Cells.Item(51, 5).Value = 551.4199829101562
$ws.Cells.Item(51, 6).Value = 240.6799926757812
$ws.Cells.Item(51, 7).Value = 216.8399963378906
$ws.Cells.Item(51, 8).Value = 168.1199951171875
$ws.Cells.Item(51, 9).Value = 385.3200073242188
$ws.Cells.Item(51, 10).Value = 559.1099853515625
$ws.Cells.Item(51, 11).Value = 248.2899932861328
$ws.Cells.Item(51, 12).Value = 208.4199981689453
$ws.Cells.Item(51, 13).Value = 164.0700073242188
$ws.Cells.Item(51, 14).Value = 377.4500122070312
$ws.Cells.Item(51, 15).Value = 549.6799926757812
$ws.Cells.Item(51, 16).Value = 232.6000061035156
$ws.Cells.Item(51, 17).Value = 215.9499969482422
$ws.Cells.Item(51, 18).Value = 167.9799957275391
$ws.Cells.Item(51, 19).Value = 383.1600036621094
$ws.Cells.Item(51, 20).Value = 558.489990234375
$ws.Cells.Item(51, 21).Value = 248.1300048828125
$ws.Cells.Item(51, 22).Value = 61368300
$ws.Cells.Item(51, 23).Value = 15206200
$ws.Cells.Item(51, 24).Value = 20473000
$ws.Cells.Item(51, 25).Value = 74079400
$ws.Cells.Item(51, 26).Value = 114813500

# New rows 52-55: copy formatting from row 51 (A col date style), then set values
$ws.Cells.Item(51,1).Copy() | Out-Null
$ws.Cells.Item(52,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(51,1).Copy() | Out-Null
$ws.Cells.Item(53,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(51,1).Copy() | Out-Null
$ws.Cells.Item(54,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(51,1).Copy() | Out-Null
$ws.Cells.Item(55,1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(52, 1).Value = 45730
$ws.Cells.Item(52, 2).Value = 213.4900054931641
$ws.Cells.Item(52, 3).Value = 167.6199951171875
$ws.Cells.Item(52, 4).Value = 388.5599975585938
$ws.Cells.Item(52, 5).Value = 562.8099975585938
$ws.Cells.Item(52, 6).Value = 249.9799957275391
$ws.Cells.Item(52, 7).Value = 213.9499969482422
$ws.Cells.Item(52, 8).Value = 168.25
$ws.Cells.Item(52, 9).Value = 390.2300109863281
$ws.Cells.Item(52, 10).Value = 563.8300170898438
$ws.Cells.Item(52, 11).Value = 251.5800018310547
$ws.Cells.Item(52, 12).Value = 209.5800018310547
$ws.Cells.Item(52, 13).Value = 164.5099945068359
$ws.Cells.Item(52, 14).Value = 379.510009765625
$ws.Cells.Item(52, 15).Value = 551.489990234375
$ws.Cells.Item(52, 16).Value = 240.7299957275391
$ws.Cells.Item(52, 17).Value = 211.25
$ws.Cells.Item(52, 18).Value = 165.3150024414062
$ws.Cells.Item(52, 19).Value = 379.7799987792969
$ws.Cells.Item(52, 20).Value = 556.1099853515625
$ws.Cells.Item(52, 21).Value = 247.3099975585938
$ws.Cells.Item(52, 22).Value = 60107600
$ws.Cells.Item(52, 23).Value = 18611100
$ws.Cells.Item(52, 24).Value = 19952800
$ws.Cells.Item(52, 25).Value = 62660300
$ws.Cells.Item(52, 26).Value = 100242300
$ws.Cells.Item(53, 1).Value = 45733
$ws.Cells.Item(53, 2).Value = 214
$ws.Cells.Item(53, 3).Value = 166.5700073242188
$ws.Cells.Item(53, 4).Value = 388.7000122070312
$ws.Cells.Item(53, 5).Value = 567.1500244140625
$ws.Cells.Item(53, 6).Value = 238.0099945068359
$ws.Cells.Item(53, 7).Value = 215.2200012207031
$ws.Cells.Item(53, 8).Value = 168.4600067138672
$ws.Cells.Item(53, 9).Value = 392.7099914550781
$ws.Cells.Item(53, 10).Value = 569.7100219726562
$ws.Cells.Item(53, 11).Value = 245.3999938964844
$ws.Cells.Item(53, 12).Value = 209.9700012207031
$ws.Cells.Item(53, 13).Value = 165.8099975585938
$ws.Cells.Item(53, 14).Value = 385.5700073242188
$ws.Cells.Item(53, 15).Value = 562.3499755859375
$ws.Cells.Item(53, 16).Value = 232.8000030517578
$ws.Cells.Item(53, 17).Value = 213.3099975585938
$ws.Cells.Item(53, 18).Value = 167.3249969482422
$ws.Cells.Item(53, 19).Value = 386.7000122070312
$ws.Cells.Item(53, 20).Value = 562.7899780273438
$ws.Cells.Item(53, 21).Value = 245.0599975585938
$ws.Cells.Item(53, 22).Value = 48073400
$ws.Cells.Item(53, 23).Value = 17839100
$ws.Cells.Item(53, 24).Value = 22474300
$ws.Cells.Item(53, 25).Value = 49008700
$ws.Cells.Item(53, 26).Value = 111900600
$ws.Cells.Item(54, 1).Value = 45734
$ws.Cells.Item(54, 2).Value = 212.6900024414062
$ws.Cells.Item(54, 3).Value = 162.6699981689453
$ws.Cells.Item(54, 4).Value = 383.5199890136719
$ws.Cells.Item(54, 5).Value = 561.02001953125
$ws.Cells.Item(54, 6).Value = 225.3099975585938
$ws.Cells.Item(54, 7).Value = 215.1499938964844
$ws.Cells.Item(54, 8).Value = 166.4400024414062
$ws.Cells.Item(54, 9).Value = 387.3699951171875
$ws.Cells.Item(54, 10).Value = 565.02001953125
$ws.Cells.Item(54, 11).Value = 230.1000061035156
$ws.Cells.Item(54, 12).Value = 211.4900054931641
$ws.Cells.Item(54, 13).Value = 158.8000030517578
$ws.Cells.Item(54, 14).Value = 381.1000061035156
$ws.Cells.Item(54, 15).Value = 559.0599975585938
$ws.Cells.Item(54, 16).Value = 222.2799987792969
$ws.Cells.Item(54, 17).Value = 214.1600036621094
$ws.Cells.Item(54, 18).Value = 165.9600067138672
$ws.Cells.Item(54, 19).Value = 387.0700073242188
$ws.Cells.Item(54, 20).Value = 564.7999877929688
$ws.Cells.Item(54, 21).Value = 228.1600036621094
$ws.Cells.Item(54, 22).Value = 42432400
$ws.Cells.Item(54, 23).Value = 24616800
$ws.Cells.Item(54, 24).Value = 19486900
$ws.Cells.Item(54, 25).Value = 66041400
$ws.Cells.Item(54, 26).Value = 111477600
$ws.Cells.Item(55, 1).Value = 45735
$ws.Cells.Item(55, 2).Value = 215.2400054931641
$ws.Cells.Item(55, 3).Value = 166.2799987792969
$ws.Cells.Item(55, 4).Value = 387.8200073242188
$ws.Cells.Item(55, 5).Value = 567.1300048828125
$ws.Cells.Item(55, 6).Value = 235.8600006103516
$ws.Cells.Item(55, 7).Value = 218.7599945068359
$ws.Cells.Item(55, 8).Value = 168.1329956054688
$ws.Cells.Item(55, 9).Value = 389.6799926757812
$ws.Cells.Item(55, 10).Value = 570.9500122070312
$ws.Cells.Item(55, 11).Value = 241.4100036621094
$ws.Cells.Item(55, 12).Value = 213.75
$ws.Cells.Item(55, 13).Value = 163.0500030517578
$ws.Cells.Item(55, 14).Value = 384
$ws.Cells.Item(55, 15).Value = 561.6300048828125
$ws.Cells.Item(55, 16).Value = 229.1999969482422
$ws.Cells.Item(55, 17).Value = 214.2200012207031
$ws.Cells.Item(55, 18).Value = 163.9149932861328
$ws.Cells.Item(55, 19).Value = 385.5299987792969
$ws.Cells.Item(55, 20).Value = 562.8300170898438
$ws.Cells.Item(55, 21).Value = 231.6100006103516
$ws.Cells.Item(55, 22).Value = 54336700
$ws.Cells.Item(55, 23).Value = 24932500
$ws.Cells.Item(55, 24).Value = 19165300
$ws.Cells.Item(55, 25).Value = 66392000
$ws.Cells.Item(55, 26).Value = 111250900
